$d = $word.ActiveDocument

# --- Part 1: merge "Project #" + "2" runs into a single run "Project #2" ---
$d.Content.Find.Execute("Project #2", $false, $false, $false, $false, $false, $true, 1, $false, "Project #2", 2) | Out-Null

# --- Part 2: merge the two Gulp-paragraph runs (around the removed bookmark) into one run ---
$gulpText = "It was easy to follow along the Gulp tutorial video, but I had trouble passing the part about sass. Mine was not working for some reason but I had taken the time to watch the video again to retrace my steps."
$d.Content.Find.Execute($gulpText, $false, $false, $false, $false, $false, $true, 1, $false, $gulpText, 2) | Out-Null

# --- Part 3: insert two new paragraphs after the "pixelated" bullet, before the trailing empty paragraph ---
$findRange = $d.Content
$findRange.Find.Execute("pixelated", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$findRange.Expand(4) | Out-Null   # wdParagraph -> expand to the whole paragraph

$insertAt = $findRange.End - 1    # just before that paragraph's own paragraph mark
$insertPoint = $d.Range($insertAt, $insertAt)

$newParasXml = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:rFonts w:ascii="Lucida Sans" w:hAnsi="Lucida Sans"/><w:sz w:val="28"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Lucida Sans" w:hAnsi="Lucida Sans"/><w:sz w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Lucida Sans" w:hAnsi="Lucida Sans"/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve">I was able to get my font working for majority of the project until I had to style it around using the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Lucida Sans" w:hAnsi="Lucida Sans"/><w:sz w:val="28"/></w:rPr><w:t>scss</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Lucida Sans" w:hAnsi="Lucida Sans"/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve"> format and then it stopped working on me. I want to say I almost tried doing everything I knew to fix it and asking my peers but the font still was not working.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

$flatOpc = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + "<w:body>$newParasXml</w:body>" + '</w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertPoint.InsertXML($flatOpc)
